# Update the cryptocurrency price/volume snapshot (Price column D, Volume(1h) column E).
# Numeric-looking Price values are written with a leading apostrophe so Excel keeps
# them as text (matching the original inlineStr cell type) instead of coercing them
# to numbers and dropping significant trailing zeros (e.g. "139.90" -> 139.9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.647.15"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "3.134.86"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'528.52"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "'139.90"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.136.60"
$ws.Range("E8").Value = "  +0.74%  "
$ws.Range("D9").Value = "'0.443"
$ws.Range("E9").Value = "  +2.29%  "
$ws.Range("D10").Value = "'7.18"
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").Value = "'0.108"
$ws.Range("E11").Value = "  -1.41%  "
$ws.Range("D12").Value = "'0.396"
$ws.Range("E12").Value = "  +2.89%  "
$ws.Range("D13").Value = "3.687.37"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("D15").Value = "'25.52"
$ws.Range("E15").Value = "  -2.92%  "
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "57.819.50"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "3.147.69"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").Value = "'6.09"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").Value = "'12.82"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "'7.93"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").Value = "'354.72"
$ws.Range("E22").Value = "  +5.16%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "'68.55"
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("D25").Value = "'0.508"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").Value = "'0.170"
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").Value = "0.0₃0928"
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("D29").Value = "'7.42"
$ws.Range("E29").Value = "  +2.40%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "'6.34"
$ws.Range("E31").Value = "  -3.70%  "
$ws.Range("D32").Value = "'1.89"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("D33").Value = "'21.13"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").Value = "'1.19"
$ws.Range("E34").Value = "  -1.83%  "
$ws.Range("D35").Value = "'4.91"
$ws.Range("E35").Value = "  +5.53%  "
$ws.Range("D36").Value = "'157.62"
$ws.Range("E36").Value = "  +2.34%  "
$ws.Range("D37").Value = "'6.17"
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("D38").Value = "'26.02"
$ws.Range("D39").Value = "'1.27"
$ws.Range("E39").Value = "  -2.85%  "
$ws.Range("D40").Value = "'0.0669"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").Value = "'1.62"
$ws.Range("E41").Value = "  +7.93%  "
$ws.Range("D42").Value = "'4.11"
$ws.Range("E42").Value = "  +4.96%  "
$ws.Range("D43").Value = "'0.704"
$ws.Range("E43").Value = "  +2.54%  "
$ws.Range("D44").Value = "3.185.36"
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("D45").Value = "'0.0272"
$ws.Range("E45").Value = "  +4.75%  "
$ws.Range("D46").Value = "'36.62"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").Value = "2.329.39"
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("D49").Value = "'0.985"
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("D50").Value = "'6.06"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").Value = "'20.25"
$ws.Range("E51").Value = "  -3.31%  "
